$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# "About" sheet: insert a two-row gap before the old row 18 and add a
# new note about doubling the BLUE Shifts potential.
# ---------------------------------------------------------------------
$about = $wb.Worksheets.Item("About")
$about.Rows("18:19").Insert()
$about.Range("A18").Value = "We allow for twice the potential identified in the BLUE Shifts scenario."

# ---------------------------------------------------------------------
# "PCiCDTdtTDM" sheet: double every lever value coming from Calcs.
# ---------------------------------------------------------------------
$pci = $wb.Worksheets.Item("PCiCDTdtTDM")
$pci.Range("B2").Formula = "=Calcs!B5*2"
$pci.Range("B3").Formula = "=Calcs!C5*2"
$pci.Range("C3").Formula = "=Calcs!B11*2"
$pci.Range("B4").Formula = "=Calcs!D5*2"
$pci.Range("B5").Formula = "=Calcs!E5*2"
$pci.Range("C5").Formula = "=Calcs!C11*2"
$pci.Range("B6").Formula = "=Calcs!F5*2"
$pci.Range("B7").Formula = "=Calcs!G5*2"
$pci.Range("C6").Select() | Out-Null

# Restore "About" as the active sheet/selection, matching the source edit.
$about.Activate() | Out-Null
$about.Rows("19:19").EntireRow.Select() | Out-Null
